$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.835141539573669
$ws.Range("B1").Value = 3.202824354171753
$ws.Range("C1").Value = 2.860763072967529
$ws.Range("D1").Value = 1.646771669387817
$ws.Range("E1").Value = 0.9598884582519531
